$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "1706"
$ws.Range("E17").Value = "1705"
$ws.Range("E18").Value = "1704"
$ws.Range("E19").Value = "1703"
$ws.Range("E20").Value = "1702"
$ws.Range("E21").Value = "1701"
